# PROGI_2023_Posterized_prezentacija.pptx - "rad na draftu prezentacije"
#
# Slide 4 ("Opis zadatka" / Content Placeholder 2): rewrite the bullet
# list describing the project idea.
# Slide 6 ("Pregled nefunkc. zahtjeva" / Content Placeholder 2): fill in
# the previously-empty bullet list of non-functional requirements.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 4 - "Opis zadatka"
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$tr4 = $sh4.TextFrame.TextRange

$lines4 = @(
    "razvijanje web aplikacije za stručne konferencije",
    "mogućnost pregledavanja i glasanja za radove u realnom vremenu",
    "intuitivno korisničko iskustvo prilagođeno potrebama sudionika konferencija",
    "integracija s umjetnom inteligencijom za personalizirane preporuke, društvene mreže, virtualne stvarnosti za obilazak dvorane,…",
    "slične platforme su Whova, EventMobi i Attendify"
)

$tr4.Text = [string]::Join([char]13, $lines4)
$tr4.Font.Size = 20

Write-Host "Slide4 paragraph count:" $tr4.Paragraphs().Count
